# Korean - 한국어(대한민국) Config Added
# Translate the English report-template headers (row 1) to Korean,
# for both the worksheet cells and the backing table ("테ーブル1") columns,
# which Excel keeps in sync automatically when the header cell text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "워크플로우 이름"   # Workflow Filename
$ws.Range("B1").Value = "내부 경로"         # Internal Path
$ws.Range("C1").Value = "대상"              # Target
$ws.Range("D1").Value = "이슈"              # Issue
$ws.Range("E1").Value = "조치"              # Action
$ws.Range("F1").Value = "제안"              # Suggestion

# Move the active selection to A2, matching the saved cursor position.
$ws.Range("A2").Select() | Out-Null
